$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Radhames Capellan" (was at C4) and "Goretti Guerrero" (was at B10) already
# finished the get/remove algorithm, so pull them out of the pair grid. Grab
# their existing cell formatting before clearing the cells.
$ws.Range("C4").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("C4").Clear()
$ws.Range("B10").Clear()

# ...and append them to the "Individual" roster in column E, continuing the
# existing list that runs from E8 down to E23.
$ws.Range("E24").Value = "Radhames Capellan"
$ws.Range("E25").Value = "Goretti Guerrero"

# Brand new team member added for the same exercise, highlighted with a
# left/right box border like the rest of the roster.
$ws.Range("E26").Value = "Marlene Villalobos"
$ws.Range("E26").Borders.Item(7).Weight = 2
$ws.Range("E26").Borders.Item(10).Weight = 2

# Leave the selection where the user would land after typing the last entry.
$ws.Range("E27").Select()
